$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Sema5a"
$ws.Cells.Item(2,3).Value2 = "Met"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 2
$ws.Cells.Item(2,6).Value2 = 0.6666666666666666
$ws.Cells.Item(2,7).Value2 = 0.4290636666666667
$ws.Cells.Item(2,8).Value2 = 1.287191
$ws.Cells.Item(2,9).Value2 = 0.0126431569814401
$ws.Cells.Item(2,10).Value2 = 0.0126431569814401
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 1.847798333333333
$ws.Cells.Item(2,14).Value2 = 5.543395
$ws.Cells.Item(2,15).Value2 = 0.05039680725746681
$ws.Cells.Item(2,16).Value2 = 0.05039680725746681
$ws.Cells.Item(2,17).Value2 = 0.7928231281605557
$ws.Cells.Item(2,18).Value2 = 7.135408153445001
$ws.Cells.Item(2,19).Value2 = 0.0006371747455195324
$ws.Cells.Item(2,20).Value2 = 0.0006371747455195323

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Sema5a"
$ws.Cells.Item(3,3).Value2 = "Met"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 2
$ws.Cells.Item(3,6).Value2 = 0.6666666666666666
$ws.Cells.Item(3,7).Value2 = 0.4290636666666667
$ws.Cells.Item(3,8).Value2 = 1.287191
$ws.Cells.Item(3,9).Value2 = 0.0126431569814401
$ws.Cells.Item(3,10).Value2 = 0.0126431569814401
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 0.4798556666666667
$ws.Cells.Item(3,14).Value2 = 1.439567
$ws.Cells.Item(3,15).Value2 = 0.0130875719001099
$ws.Cells.Item(3,16).Value2 = 0.0130875719001099
$ws.Cells.Item(3,17).Value2 = 0.2058886318107778
$ws.Cells.Item(3,18).Value2 = 1.852997686297
$ws.Cells.Item(3,19).Value2 = 0.0001654682260389737
$ws.Cells.Item(3,20).Value2 = 0.0001654682260389737

# Row 4: ECs -> sCs
$ws.Cells.Item(4,1).Value2 = "ECs"
$ws.Cells.Item(4,2).Value2 = "Sema5a"
$ws.Cells.Item(4,3).Value2 = "Met"
$ws.Cells.Item(4,4).Value2 = "sCs"
$ws.Cells.Item(4,5).Value2 = 2
$ws.Cells.Item(4,6).Value2 = 0.6666666666666666
$ws.Cells.Item(4,7).Value2 = 0.4290636666666667
$ws.Cells.Item(4,8).Value2 = 1.287191
$ws.Cells.Item(4,9).Value2 = 0.0126431569814401
$ws.Cells.Item(4,10).Value2 = 0.0126431569814401
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 34.337334
$ws.Cells.Item(4,14).Value2 = 103.012002
$ws.Cells.Item(4,15).Value2 = 0.9365156208424232
$ws.Cells.Item(4,16).Value2 = 0.9365156208424232
$ws.Cells.Item(4,17).Value2 = 14.732902429598
$ws.Cells.Item(4,18).Value2 = 132.596121866382
$ws.Cells.Item(4,19).Value2 = 0.01184051400988159
$ws.Cells.Item(4,20).Value2 = 0.01184051400988159

# Row 5: FAPs -> ECs
$ws.Cells.Item(5,1).Value2 = "FAPs"
$ws.Cells.Item(5,2).Value2 = "Sema5a"
$ws.Cells.Item(5,3).Value2 = "Met"
$ws.Cells.Item(5,4).Value2 = "ECs"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 26.436942
$ws.Cells.Item(5,8).Value2 = 79.31082599999999
$ws.Cells.Item(5,9).Value2 = 0.7790135445677298
$ws.Cells.Item(5,10).Value2 = 0.7790135445677296
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 1.847798333333333
$ws.Cells.Item(5,14).Value2 = 5.543395
$ws.Cells.Item(5,15).Value2 = 0.05039680725746681
$ws.Cells.Item(5,16).Value2 = 0.05039680725746681
$ws.Cells.Item(5,17).Value2 = 48.85013736603
$ws.Cells.Item(5,18).Value2 = 439.6512362942699
$ws.Cells.Item(5,19).Value2 = 0.03925979545653591
$ws.Cells.Item(5,20).Value2 = 0.0392597954565359

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6,1).Value2 = "FAPs"
$ws.Cells.Item(6,2).Value2 = "Sema5a"
$ws.Cells.Item(6,3).Value2 = "Met"
$ws.Cells.Item(6,4).Value2 = "FAPs"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 26.436942
$ws.Cells.Item(6,8).Value2 = 79.31082599999999
$ws.Cells.Item(6,9).Value2 = 0.7790135445677298
$ws.Cells.Item(6,10).Value2 = 0.7790135445677296
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 0.4798556666666667
$ws.Cells.Item(6,14).Value2 = 1.439567
$ws.Cells.Item(6,15).Value2 = 0.0130875719001099
$ws.Cells.Item(6,16).Value2 = 0.0130875719001099
$ws.Cells.Item(6,17).Value2 = 12.685916428038
$ws.Cells.Item(6,18).Value2 = 114.173247852342
$ws.Cells.Item(6,19).Value2 = 0.01019539577568963
$ws.Cells.Item(6,20).Value2 = 0.01019539577568963

# Row 7: FAPs -> sCs
$ws.Cells.Item(7,1).Value2 = "FAPs"
$ws.Cells.Item(7,2).Value2 = "Sema5a"
$ws.Cells.Item(7,3).Value2 = "Met"
$ws.Cells.Item(7,4).Value2 = "sCs"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 26.436942
$ws.Cells.Item(7,8).Value2 = 79.31082599999999
$ws.Cells.Item(7,9).Value2 = 0.7790135445677298
$ws.Cells.Item(7,10).Value2 = 0.7790135445677296
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 34.337334
$ws.Cells.Item(7,14).Value2 = 103.012002
$ws.Cells.Item(7,15).Value2 = 0.9365156208424232
$ws.Cells.Item(7,16).Value2 = 0.9365156208424232
$ws.Cells.Item(7,17).Value2 = 907.7741073926279
$ws.Cells.Item(7,18).Value2 = 8169.966966533651
$ws.Cells.Item(7,19).Value2 = 0.7295583533355042
$ws.Cells.Item(7,20).Value2 = 0.7295583533355041

# Row 8: sCs -> ECs
$ws.Cells.Item(8,1).Value2 = "sCs"
$ws.Cells.Item(8,2).Value2 = "Sema5a"
$ws.Cells.Item(8,3).Value2 = "Met"
$ws.Cells.Item(8,4).Value2 = "ECs"
$ws.Cells.Item(8,5).Value2 = 3
$ws.Cells.Item(8,6).Value2 = 1
$ws.Cells.Item(8,7).Value2 = 7.070428666666666
$ws.Cells.Item(8,8).Value2 = 21.211286
$ws.Cells.Item(8,9).Value2 = 0.2083432984508302
$ws.Cells.Item(8,10).Value2 = 0.2083432984508302
$ws.Cells.Item(8,11).Value2 = 3
$ws.Cells.Item(8,12).Value2 = 1
$ws.Cells.Item(8,13).Value2 = 1.847798333333333
$ws.Cells.Item(8,14).Value2 = 5.543395
$ws.Cells.Item(8,15).Value2 = 0.05039680725746681
$ws.Cells.Item(8,16).Value2 = 0.05039680725746681
$ws.Cells.Item(8,17).Value2 = 13.06472630621889
$ws.Cells.Item(8,18).Value2 = 117.58253675597
$ws.Cells.Item(8,19).Value2 = 0.01049983705541137
$ws.Cells.Item(8,20).Value2 = 0.01049983705541137

# Row 9: sCs -> FAPs
$ws.Cells.Item(9,1).Value2 = "sCs"
$ws.Cells.Item(9,2).Value2 = "Sema5a"
$ws.Cells.Item(9,3).Value2 = "Met"
$ws.Cells.Item(9,4).Value2 = "FAPs"
$ws.Cells.Item(9,5).Value2 = 3
$ws.Cells.Item(9,6).Value2 = 1
$ws.Cells.Item(9,7).Value2 = 7.070428666666666
$ws.Cells.Item(9,8).Value2 = 21.211286
$ws.Cells.Item(9,9).Value2 = 0.2083432984508302
$ws.Cells.Item(9,10).Value2 = 0.2083432984508302
$ws.Cells.Item(9,11).Value2 = 3
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 0.4798556666666667
$ws.Cells.Item(9,14).Value2 = 1.439567
$ws.Cells.Item(9,15).Value2 = 0.0130875719001099
$ws.Cells.Item(9,16).Value2 = 0.0130875719001099
$ws.Cells.Item(9,17).Value2 = 3.392785261462444
$ws.Cells.Item(9,18).Value2 = 30.535067353162
$ws.Cells.Item(9,19).Value2 = 0.002726707898381295
$ws.Cells.Item(9,20).Value2 = 0.002726707898381295

# Row 10: sCs -> sCs
$ws.Cells.Item(10,1).Value2 = "sCs"
$ws.Cells.Item(10,2).Value2 = "Sema5a"
$ws.Cells.Item(10,3).Value2 = "Met"
$ws.Cells.Item(10,4).Value2 = "sCs"
$ws.Cells.Item(10,5).Value2 = 3
$ws.Cells.Item(10,6).Value2 = 1
$ws.Cells.Item(10,7).Value2 = 7.070428666666666
$ws.Cells.Item(10,8).Value2 = 21.211286
$ws.Cells.Item(10,9).Value2 = 0.2083432984508302
$ws.Cells.Item(10,10).Value2 = 0.2083432984508302
$ws.Cells.Item(10,11).Value2 = 3
$ws.Cells.Item(10,12).Value2 = 1
$ws.Cells.Item(10,13).Value2 = 34.337334
$ws.Cells.Item(10,14).Value2 = 103.012002
$ws.Cells.Item(10,15).Value2 = 0.9365156208424232
$ws.Cells.Item(10,16).Value2 = 0.9365156208424232
$ws.Cells.Item(10,17).Value2 = 242.7796706505079
$ws.Cells.Item(10,18).Value2 = 2185.017035854572
$ws.Cells.Item(10,19).Value2 = 0.1951167534970375
$ws.Cells.Item(10,20).Value2 = 0.1951167534970375
